$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.442.41'
$ws.Range("E2").Value = '  -0.20%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.520.36'
$ws.Range("E3").Value = '  -2.22%  '

# Row 4
$ws.Range("E4").Value = '  +0.06%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.69'
$ws.Range("E5").Value = '  -0.63%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.16'
$ws.Range("E6").Value = '  -3.94%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.519.44'
$ws.Range("E7").Value = '  -2.19%  '

# Row 8
$ws.Range("E8").Value = '  -0.29%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.515'
$ws.Range("E9").Value = '  +5.45%  '

# Row 10
$ws.Range("B10").Value = 'Toncoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.71'
$ws.Range("E10").Value = '  -4.63%  '

# Row 11
$ws.Range("B11").Value = 'Dogecoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.130'
$ws.Range("E11").Value = '  -4.34%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.408'
$ws.Range("E12").Value = '  -1.88%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.126.25'
$ws.Range("E13").Value = '  -2.09%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000194'
$ws.Range("E14").Value = '  -7.05%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.60'
$ws.Range("E15").Value = '  -3.93%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.523.58'
$ws.Range("E16").Value = '  -2.54%  '

# Row 17
$ws.Range("E17").Value = '  +0.57%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.373.57'
$ws.Range("E18").Value = '  -0.46%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.86'
$ws.Range("E19").Value = '  -5.97%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.16'
$ws.Range("E20").Value = '  -3.61%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.59'
$ws.Range("E21").Value = '  -3.50%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '421.90'
$ws.Range("E22").Value = '  -1.50%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.588'
$ws.Range("E23").Value = '  -5.06%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.19'
$ws.Range("E24").Value = '  -2.22%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.669.70'
$ws.Range("E25").Value = '  -2.20%  '

# Row 26
$ws.Range("E26").Value = '  +0.00%  '

# Row 27
$ws.Range("E27").Value = '  -7.13%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.89'
$ws.Range("E28").Value = '  -5.03%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.46'
$ws.Range("E29").Value = '  -2.60%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.92'
$ws.Range("E30").Value = '  -5.68%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.00'
$ws.Range("E31").Value = '  +0.00%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.534.91'
$ws.Range("E32").Value = '  -1.79%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.155'
$ws.Range("E33").Value = '  -2.00%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '24.15'
$ws.Range("E34").Value = '  -5.28%  '

# Row 35
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.34'
$ws.Range("E36").Value = '  -9.19%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.52'
$ws.Range("E37").Value = '  -4.65%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.62'
$ws.Range("E38").Value = '  -4.67%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '174.20'
$ws.Range("E39").Value = '  -2.15%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.20'
$ws.Range("E40").Value = '  -8.23%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0817'
$ws.Range("E41").Value = '  -4.64%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.97'
$ws.Range("E42").Value = '  -5.39%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.857'
$ws.Range("E43").Value = '  -4.79%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.34'
$ws.Range("E44").Value = '  -1.45%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.76'
$ws.Range("E45").Value = '  -7.84%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").Value = '  +0.16%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.35'
$ws.Range("E47").Value = '  -8.12%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.06'
$ws.Range("E48").Value = '  -1.94%  '

# Row 49
$ws.Range("E49").Value = '  -6.34%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.75'
$ws.Range("E50").Value = '  -5.60%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.906'
$ws.Range("E51").Value = '  -5.25%  '
